# #5: property boat&car done
# Populate the "汽車" (car) sheet (3rd worksheet) with the full set of
# columns used by the other property sheets (name/capacity/owner/...)
# instead of the old partial land-style row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# --- First extend the formatting (style) of the existing header/data
#     rows across the new columns H:N, matching the s=1 (header) /
#     s=2 (data) styles already used in columns B:G.
$ws.Range("B1:G1").Copy() | Out-Null
$ws.Range("H1:N1").PasteSpecial(-4122) | Out-Null

$ws.Range("B2:G2").Copy() | Out-Null
$ws.Range("H2:N2").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# --- Header row (row 1): now the generic property-schema header,
#     same as the 土地/建物 sheets, with "capacity" swapped in for
#     "area" since this sheet tracks vehicles.
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "capacity"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "register_date"
$ws.Range("F1").Value = "register_reason"
$ws.Range("G1").Value = "acquire_value"
$ws.Range("H1").Value = "property_category"
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# --- Data row (row 2): fill in the newly-added trailing columns;
#     existing B2:G2 values are unchanged.
$ws.Range("H2").Value = "land"
$ws.Range("I2").Value = "normal"

# J2 ("2012-04-30") must stay literal text, not be reinterpreted as a
# date serial -- format the cell as Text first.
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "2012-04-30"

$ws.Range("K2").Value = "趙天麟"
$ws.Range("L2").Value = 1761
$ws.Range("M2").Value = "tmp58581"
$ws.Range("N2").Value = 31
